# Weekly fruit/vegetable price update: a new weekly record is inserted
# into the "Choclo" (corn) price history at row 27, pushing the
# subsequent historical rows (27-54) down by one (to 28-55).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 27; existing rows 27-54 shift down to 28-55
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly price record
$ws.Cells.Item(27,1).Value = 11
$ws.Cells.Item(27,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(27,3).Value = "Bíobío"
$ws.Cells.Item(27,4).Value = 44483
$ws.Cells.Item(27,5).Value = 8
$ws.Cells.Item(27,6).Value = 100112024
$ws.Cells.Item(27,7).Value = "Choclo"
$ws.Cells.Item(27,8).Value = "Dulce o Americano"
$ws.Cells.Item(27,9).Value = "Primera"
$ws.Cells.Item(27,10).Value = 450
$ws.Cells.Item(27,11).Value = 25000
$ws.Cells.Item(27,12).Value = 26000
$ws.Cells.Item(27,13).Value = 25556
$ws.Cells.Item(27,14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(27,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27,16).Value = 365
$ws.Cells.Item(27,17).Value = 70
$ws.Cells.Item(27,18).Value = "Hortaliza"
